$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.506.80'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.02'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.05'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3522'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.237'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07750'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9990'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.59'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +12.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.635'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.228'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.811.77'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001128'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06756'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '86.98'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +7.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9985'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.84'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +10.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.533'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.20'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.520.52'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.463'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.686'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +15.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.505'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +16.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.21'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.014.61'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '136.89'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +7.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.395'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +7.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.075'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.95'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08795'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.720'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.675'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.63%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +16.22%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06573'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.83%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02428'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.13%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2272'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.020'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.294'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.05'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6626'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9982'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.047'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.36%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.20'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07363'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.83'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.81%  '
